# This script populates the "O:X" (Skip-gram) result columns for the
# Word2Vec worksheet's 10/50/100/150/200-feature blocks (rows 135-139,
# 143-147, 151-153), mirroring the data added in the tracked commit
# ("added sections needed for the report"). A handful of cells in the
# original workbook already carry slightly different styling (picked up
# from earlier ad-hoc data entry/highlighting), so those are replicated
# explicitly below instead of being left to whatever the default style
# resolution would produce.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Word2Vec")
$ws.Activate()

$styleDonor = $ws.Range("O63")

$ws.Range("O135").Value = 0.87375000000000003
$ws.Range("P135").Value = 0.87250000000000005
$ws.Range("Q135").Value = 0.84875
$ws.Range("R135").Value = 0.88249999999999995
$ws.Range("S135").Value = 0.86624999999999996
$ws.Range("T135").Value = 0.86875000000000002
$ws.Range("U135").Value = 0.86750000000000005
$ws.Range("V135").Value = 0.87250000000000005
$ws.Range("W135").Value = 0.84750000000000003
$ws.Range("X135").Value = 0.87124999999999997
$ws.Range("O136").Value = 0.88
$ws.Range("P136").Value = 0.88624999999999998
$ws.Range("Q136").Value = 0.87749999999999995
$ws.Range("R136").Value = 0.88
$ws.Range("S136").Value = 0.88
$ws.Range("T136").Value = 0.88249999999999995
$ws.Range("U136").Value = 0.87250000000000005
$ws.Range("V136").Value = 0.89875000000000005
$ws.Range("W136").Value = 0.87250000000000005
$styleDonor.Copy()
$ws.Range("X136").PasteSpecial(-4122)
$ws.Range("X136").Value = 0.88124999999999998
$ws.Range("O137").Value = 0.87875000000000003
$ws.Range("P137").Value = 0.88624999999999998
$ws.Range("Q137").Value = 0.86624999999999996
$ws.Range("R137").Value = 0.88
$ws.Range("S137").Value = 0.88624999999999998
$ws.Range("T137").Value = 0.88124999999999998
$ws.Range("U137").Value = 0.87624999999999997
$ws.Range("V137").Value = 0.89875000000000005
$ws.Range("W137").Value = 0.86499999999999999
$ws.Range("X137").Value = 0.88749999999999996
$ws.Range("O138").Value = 0.88
$ws.Range("P138").Value = 0.88500000000000001
$ws.Range("Q138").Value = 0.86499999999999999
$ws.Range("R138").Value = 0.88624999999999998
$ws.Range("S138").Value = 0.88124999999999998
$ws.Range("T138").Value = 0.88249999999999995
$ws.Range("U138").Value = 0.87875000000000003
$ws.Range("V138").Value = 0.89749999999999996
$ws.Range("W138").Value = 0.86875000000000002
$ws.Range("X138").Value = 0.88624999999999998
$ws.Range("O139").Value = 0.88
$ws.Range("P139").Value = 0.88375000000000004
$ws.Range("Q139").Value = 0.87
$ws.Range("R139").Value = 0.88749999999999996
$ws.Range("S139").Value = 0.88624999999999998
$ws.Range("T139").Value = 0.88249999999999995
$ws.Range("U139").Value = 0.87875000000000003
$ws.Range("V139").Value = 0.89875000000000005
$ws.Range("W139").Value = 0.87
$ws.Range("X139").Value = 0.88375000000000004
$ws.Range("O143").Value = 0.74375000000000002
$ws.Range("P143").Value = 0.7
$ws.Range("Q143").Value = 0.70625000000000004
$ws.Range("R143").Value = 0.73124999999999996
$ws.Range("S143").Value = 0.78125
$ws.Range("T143").Value = 0.75624999999999998
$ws.Range("U143").Value = 0.69374999999999998
$ws.Range("V143").Value = 0.7
$ws.Range("W143").Value = 0.76249999999999996
$ws.Range("X143").Value = 0.73124999999999996
$ws.Range("O144").Value = 0.76875000000000004
$ws.Range("P144").Value = 0.77500000000000002
$ws.Range("Q144").Value = 0.76249999999999996
$ws.Range("R144").Value = 0.78125
$ws.Range("S144").Value = 0.79374999999999996
$ws.Range("T144").Value = 0.74375000000000002
$ws.Range("U144").Value = 0.74375000000000002
$ws.Range("V144").Value = 0.76249999999999996
$ws.Range("W144").Value = 0.85
$ws.Range("X144").Value = 0.78749999999999998
$ws.Range("O145").Value = 0.76875000000000004
$ws.Range("P145").Value = 0.76875000000000004
$ws.Range("Q145").Value = 0.71250000000000002
$ws.Range("R145").Value = 0.8
$ws.Range("S145").Value = 0.80625000000000002
$ws.Range("T145").Value = 0.79374999999999996
$ws.Range("U145").Value = 0.73750000000000004
$ws.Range("V145").Value = 0.76875000000000004
$ws.Range("W145").Value = 0.82499999999999996
$ws.Range("X145").Value = 0.83125000000000004
$ws.Range("O146").Value = 0.78749999999999998
$ws.Range("P146").Value = 0.75
$ws.Range("Q146").Value = 0.73124999999999996
$ws.Range("R146").Value = 0.78125
$ws.Range("S146").Value = 0.77500000000000002
$ws.Range("T146").Value = 0.79374999999999996
$ws.Range("U146").Value = 0.71250000000000002
$ws.Range("V146").Value = 0.8
$ws.Range("W146").Value = 0.84375
$ws.Range("X146").Value = 0.83750000000000002
$ws.Range("O147").Value = 0.76249999999999996
$ws.Range("P147").Value = 0.74375000000000002
$ws.Range("Q147").Value = 0.74375000000000002
$ws.Range("R147").Value = 0.79374999999999996
$ws.Range("S147").Value = 0.80625000000000002
$ws.Range("T147").Value = 0.8125
$ws.Range("U147").Value = 0.75
$ws.Range("V147").Value = 0.78125
$ws.Range("W147").Value = 0.83750000000000002
$styleDonor.Copy()
$ws.Range("X147").PasteSpecial(-4122)
$ws.Range("X147").Value = 0.85
$ws.Range("O151").Value = 0.81440000000000001
$ws.Range("P151").Value = 0.78879999999999995
$ws.Range("Q151").Value = 0.78520000000000001
$ws.Range("R151").Value = 0.80879999999999996
$ws.Range("S151").Value = 0.80559999999999998
$ws.Range("T151").Value = 0.79320000000000002
$ws.Range("U151").Value = 0.78039999999999998
$ws.Range("V151").Value = 0.78839999999999999
$ws.Range("W151").Value = 0.82320000000000004
$ws.Range("X151").Value = 0.7944
$ws.Range("O152").Value = 0.83960000000000001
$ws.Range("P152").Style = "Normal"
$ws.Range("P152").Value = 0.82399999999999995
$ws.Range("Q152").Value = 0.81440000000000001
$ws.Range("R152").Value = 0.83079999999999998
$ws.Range("S152").Value = 0.82720000000000005
$ws.Range("T152").Value = 0.82679999999999998
$ws.Range("U152").Value = 0.81440000000000001
$ws.Range("V152").Value = 0.83840000000000003
$ws.Range("W152").Value = 0.83160000000000001
$ws.Range("X152").Value = 0.81799999999999995
$ws.Range("O153").Value = 0.84199999999999997
$ws.Range("P153").Value = 0.82640000000000002
$ws.Range("Q153").Value = 0.80879999999999996
$ws.Range("R153").Value = 0.83160000000000001
$ws.Range("S153").Value = 0.83679999999999999
$ws.Range("T153").Value = 0.83399999999999996
$ws.Range("U153").Value = 0.81759999999999999
$ws.Range("V153").Value = 0.82640000000000002
$ws.Range("W153").Value = 0.84840000000000004
$ws.Range("X153").Value = 0.82520000000000004

# Row-max highlight style donor no longer needed once the two special
# cells above have borrowed its format.

# Restore the cursor position/selection shown in the final workbook.
$ws.Range("T155").Select()
